# Fill in the teacher-specific details on the exam-duty bill form.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header block: name / designation / department / term / year -----------
$ws.Range("A3").Value = "নাম: Mr. Abdul Aziz"
$ws.Range("A4").Value = "পদবী: সহকারী অধ্যাপক"
$ws.Range("G4").Value = "৪র্থ"
$ws.Range("I4").Value = "১ম"
$ws.Range("B5").Value = "সিএসই"
$ws.Range("F5").Value = "বিভাগ :সিএসই"

# --- Amount in words, below the grand-total row -----------------------------
$ws.Range("A32").Value = "কথায়:চার লক্ষ আটচল্লিশ হাজার একশত আটত্রিশ টাকা মাত্র।"

# --- Column A is now wide enough to show the longer serial labels ----------
$ws.Columns.Item(1).ColumnWidth = 13.5

# --- Row 36 grows to fit the wrapped "amount in words" style signature row -
$ws.Rows.Item(36).RowHeight = 68.4

# --- Scroll/selection state left behind by the editing session -------------
$ws.Application.Goto($ws.Range("A7"), $false)
$ws.Range("I32").Select()
